$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.930.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.790.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +12.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.81"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.788.61"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.44"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.68%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.426.46"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.780.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.962.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.15%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.21%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +18.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "482.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.77"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.50"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.11%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.58"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.09%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.938.19"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +17.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.32"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.59"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.59"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.179"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.24"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.20%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.738.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.02"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.38%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +14.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000331"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +26.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.968"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "160.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "49.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.42"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.19%  "
